# Masses astres.xlsx update
# - Table "Tableau2" header columns get a unit suffix on a 2nd line
#   ("Masse" -> "Masse\n(Kg)", etc.), which also renames the displayed
#   header cells in row 1 of Feuil1 and applies word-wrap to those cells.
# - The three formulas that reference the "Diamètre" column via a
#   structured reference are updated to the new column name.
# - J12:J21 ("Temps de rotation sur elle-même") switch from text like
#   "655,2 heures" to plain numeric hours, since the unit now lives in
#   the header.
# - I17 and I19 (previously blank "Distance planète-satellite") are
#   filled in with their distance values.
# - Selection moves to J2:J7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$nl = [char]10

# --- Rename table header columns (this also updates the Tableau2
#     ListColumns + the row-1 header cells' shared-string text) ---
$ws.Range("B1").Value = "Masse" + $nl + "(Kg)"
$ws.Range("C1").Value = "Diamètre" + $nl + "(Km)"
$ws.Range("E1").Value = "Angle de rotation " + $nl + "(deg)"
$ws.Range("F1").Value = "vitesse rotation équateur " + $nl + "(Km/h)"
$ws.Range("G1").Value = "vitesse rotation équateur " + $nl + "(rad/ds)"
$ws.Range("H1").Value = "Distance au soleil" + $nl + "(Km)"
$ws.Range("I1").Value = "Distance planète-satellite" + $nl + "(Km)"
$ws.Range("J1").Value = "Temps de rotation sur elle-même" + $nl + "(heures)"

# --- Word-wrap the header row cells whose text now spans two lines ---
$ws.Range("B1").WrapText = $true
$ws.Range("C1").WrapText = $true
$ws.Range("E1").WrapText = $true
$ws.Range("F1").WrapText = $true
$ws.Range("G1").WrapText = $true
$ws.Range("H1").WrapText = $true
$ws.Range("I1").WrapText = $true
$ws.Range("J1").WrapText = $true

# --- Update the formulas that reference [Diamètre] via a structured
#     reference so they point at the renamed column ---
$ws.Range("F2").Formula = "=  Tableau2[[#This Row],[Diamètre" + $nl + "(Km)]]*2*PI()/(24.47*24)"
$ws.Range("F13").Formula = "=(PI()*Tableau2[[#This Row],[Diamètre" + $nl + "(Km)]])/85"
$ws.Range("F14").Formula = "=(PI()*Tableau2[[#This Row],[Diamètre" + $nl + "(Km)]])/42.5"

# --- J12:J21 "Temps de rotation sur elle-même": text "NNN,N heures" -> number ---
$ws.Range("J12").Value = 655.2
$ws.Range("J13").Value = 85.2
$ws.Range("J14").Value = 42.5
$ws.Range("J15").Value = 402.8
$ws.Range("J16").Value = 173.3
$ws.Range("J17").Value = 1903.2
$ws.Range("J18").Value = 68.4
$ws.Range("J19").Value = 108
$ws.Range("J20").Value = 141.6
$ws.Range("J21").Value = 153.6

# --- Fill in the previously blank "Distance planète-satellite" cells ---
$ws.Range("I17").Value = 1070000
$ws.Range("I19").Value = 527000

# --- Move the active selection ---
$ws.Range("J2:J7").Select()
